$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Add the new data row for 27 April 2020 (27 Nisan 2020)
$ws.Range("A47").Value = 43948
$ws.Range("B47").Value = 20143
$ws.Range("C47").Value = 2131
$ws.Range("D47").Value = 95
$ws.Range("E47").Value = 4651

# Resize the table to include the new row
$tbl = $ws.ListObjects.Item("Table3")
$tbl.Resize($ws.Range("A1:E47"))

# Update the selected cell to match the diff
$ws.Range("E46").Select()
